$wb = $excel.ActiveWorkbook

# Add data for 2025-09-25: update the 2025 (column L) figures, and a handful of
# other historical cells that were corrected, across the Citywide Totals sheet,
# the By Neighborhood summary sheet, and each affected neighborhood sheet.

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4951
$ws.Range("L3").Value = 5323
$ws.Range("B4").Value = 1718
$ws.Range("E4").Value = 2057
$ws.Range("F4").Value = 1933
$ws.Range("K4").Value = 1783
$ws.Range("L4").Value = 1301
$ws.Range("L6").Value = 4504
$ws.Range("B7").Value = 23350
$ws.Range("E7").Value = 26062
$ws.Range("F7").Value = 24126
$ws.Range("K7").Value = 27575
$ws.Range("L7").Value = 16392

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 136
$ws.Range("L4").Value = 62
$ws.Range("L8").Value = 1091
$ws.Range("L11").Value = 267
$ws.Range("L19").Value = 451
$ws.Range("L20").Value = 410
$ws.Range("L21").Value = 52
$ws.Range("L25").Value = 97
$ws.Range("L27").Value = 148
$ws.Range("L29").Value = 897
$ws.Range("L33").Value = 753
$ws.Range("L34").Value = 94
$ws.Range("L36").Value = 214
$ws.Range("L37").Value = 621
$ws.Range("L42").Value = 535
$ws.Range("L44").Value = 113
$ws.Range("L48").Value = 212
$ws.Range("L52").Value = 332
$ws.Range("L53").Value = 183
$ws.Range("L54").Value = 349
$ws.Range("L55").Value = 158
$ws.Range("L57").Value = 57
$ws.Range("B63").Value = 422
$ws.Range("E63").Value = 390
$ws.Range("F63").Value = 217
$ws.Range("K63").Value = 171
$ws.Range("L63").Value = 46
$ws.Range("L64").Value = 110
$ws.Range("L65").Value = 323
$ws.Range("L67").Value = 572
$ws.Range("L72").Value = 63
$ws.Range("L76").Value = 251
$ws.Range("L79").Value = 434
$ws.Range("L85").Value = 837
$ws.Range("L90").Value = 165
$ws.Range("L92").Value = 50
$ws.Range("L95").Value = 226
$ws.Range("L96").Value = 189
$ws.Range("L99").Value = 284
$ws.Range("B101").Value = 23350
$ws.Range("E101").Value = 26062
$ws.Range("F101").Value = 24126
$ws.Range("K101").Value = 27575
$ws.Range("L101").Value = 16392

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 57
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 189

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 81
$ws.Range("L4").Value = 21
$ws.Range("L6").Value = 63
$ws.Range("L7").Value = 267

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 70
$ws.Range("L6").Value = 65

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 250
$ws.Range("L6").Value = 177
$ws.Range("L7").Value = 837

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 106
$ws.Range("L7").Value = 332

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 63
$ws.Range("L7").Value = 183

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 317
$ws.Range("L3").Value = 372
$ws.Range("L6").Value = 282
$ws.Range("L7").Value = 1091

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 205
$ws.Range("L3").Value = 261
$ws.Range("L6").Value = 226
$ws.Range("L7").Value = 753

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 226

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 210
$ws.Range("L6").Value = 173
$ws.Range("L7").Value = 621

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 118
$ws.Range("L3").Value = 100
$ws.Range("L7").Value = 323

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 78
$ws.Range("L7").Value = 284

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 170
$ws.Range("L7").Value = 572

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 85
$ws.Range("L7").Value = 349

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L6").Value = 230
$ws.Range("L7").Value = 897

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 212

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 160
$ws.Range("L4").Value = 18
$ws.Range("L6").Value = 127
$ws.Range("L7").Value = 451

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L3").Value = 33
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 51
$ws.Range("L7").Value = 251

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 178
$ws.Range("L7").Value = 535

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 49
$ws.Range("L7").Value = 158

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 143
$ws.Range("L6").Value = 97
$ws.Range("L7").Value = 434

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 30
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L6").Value = 111
$ws.Range("L7").Value = 410

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 78
$ws.Range("L7").Value = 214

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 46
$ws.Range("L7").Value = 136

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 148

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 56
$ws.Range("L7").Value = 165

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 62
